# Auto-update draw results on excel 2025-12-18T17:45:13Z
# Appends the new Pick 4 draw-result row (row 93) to the "Results" sheet,
# mirroring every prior row: Date, Game, Phase, Result, InsertedAt.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$newRow = 93

# Columns A (date-looking "2025-12-18") and C (digit-only "251218") would
# otherwise be auto-converted by Excel into a date serial / a plain number
# when assigned through .Value, so force those two cells to Text first -
# exactly as every other row in this sheet already stores them. B, D and E
# do not round-trip as numbers/dates, so they can be set directly.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025-12-18"

$ws.Range("B" + $newRow).Value = "Pick 4"

$ws.Range("C" + $newRow).NumberFormat = "@"
$ws.Range("C" + $newRow).Value = "251218"

$ws.Range("D" + $newRow).Value = "9-4-5-1"

$ws.Range("E" + $newRow).Value = "2025-12-18T21:45:13.183+04:00"
